$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 2373.5
$ws.Range("I131").Value = 2373.5
$ws.Range("K131").Value = 7120.5
$ws.Range("M131").Value = -2080.5

$ws.Range("H137").Value = 2879.4
$ws.Range("I137").Value = 3158.8
$ws.Range("K137").Value = 9476.400000000001
$ws.Range("M137").Value = -6926.400000000001

$ws.Range("H138").Value = 3473.0708
$ws.Range("I138").Value = 2322.8
$ws.Range("J138").Value = 3764.2786
$ws.Range("K138").Value = 6968.400000000001
$ws.Range("L138").Value = 11292.8358
$ws.Range("M138").Value = -1828.400000000001
$ws.Range("N138").Value = -21572.8358

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1523.2222
$ws.Range("I2").Value = 1101.2
$ws.Range("K2").Value = 1101.2
$ws.Range("M2").Value = -988.2

$ws.Range("H32").Value = 34661484
$ws.Range("J32").Value = 6537021.5
$ws.Range("L32").Value = 6537021.5
$ws.Range("N32").Value = -6537595.5

$ws.Range("H45").Value = 4809.905
$ws.Range("I45").Value = 4560.7334
$ws.Range("J45").Value = 5432.8335
$ws.Range("K45").Value = 4560.7334
$ws.Range("L45").Value = 5432.8335
$ws.Range("M45").Value = -4183.7334
$ws.Range("N45").Value = -6186.8335

$ws.Range("H61").Value = 3262.9092
$ws.Range("I61").Value = 3026.1538
$ws.Range("J61").Value = 4142.2856
$ws.Range("K61").Value = 3026.1538
$ws.Range("L61").Value = 4142.2856
$ws.Range("M61").Value = -2814.1538
$ws.Range("N61").Value = -4566.2856

$ws.Range("H63").Value = 4970
$ws.Range("I63").Value = 2500
$ws.Range("K63").Value = 2500
$ws.Range("M63").Value = -1814

$ws.Range("H66").Value = 4970
$ws.Range("I66").Value = 2500
$ws.Range("K66").Value = 12500
$ws.Range("M66").Value = -9068

$ws.Range("H110").Value = 1344.1111
$ws.Range("I110").Value = 656.7143
$ws.Range("K110").Value = 656.7143
$ws.Range("M110").Value = 1388.2857

$ws.Range("H116").Value = 1523.2222
$ws.Range("I116").Value = 1101.2
$ws.Range("K116").Value = 1101.2
$ws.Range("M116").Value = 1192.8

$ws.Range("H136").Value = 3262.9092
$ws.Range("I136").Value = 3026.1538
$ws.Range("J136").Value = 4142.2856
$ws.Range("K136").Value = 9078.4614
$ws.Range("L136").Value = 12426.8568
$ws.Range("M136").Value = -6528.4614
$ws.Range("N136").Value = -17526.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1523.2222
$ws.Range("I3").Value = 1101.2
$ws.Range("K3").Value = 1101.2
$ws.Range("M3").Value = -987.2

$ws.Range("H105").Value = 2484.9473
$ws.Range("I105").Value = 2075.875
$ws.Range("J105").Value = 4666.6665
$ws.Range("K105").Value = 2075.875
$ws.Range("L105").Value = 4666.6665
$ws.Range("M105").Value = -328.875
$ws.Range("N105").Value = -8160.6665

$ws.Range("H130").Value = 147496
$ws.Range("J130").Value = 147496
$ws.Range("L130").Value = 147496
$ws.Range("N130").Value = -157536

$ws.Range("H134").Value = 1803962.9
$ws.Range("I134").Value = 1906760.8
$ws.Range("K134").Value = 5720282.4
$ws.Range("M134").Value = -5717747.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 200000800
$ws.Range("I16").Value = 200000800
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 200000800
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -200000513
$ws.Range("N16").Value = ""

$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").Value = ""

$ws.Range("H105").Value = 1177.4445
$ws.Range("I105").Value = 887.125
$ws.Range("K105").Value = 887.125
$ws.Range("M105").Value = 859.875

$ws.Range("H107").Value = 701
$ws.Range("J107").Value = 2998
$ws.Range("L107").Value = 2998
$ws.Range("N107").Value = -6838

$ws.Range("H113").Value = 200000800
$ws.Range("I113").Value = 200000800
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 200000800
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -199998630
$ws.Range("N113").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 500
$ws.Range("I21").Value = 500
$ws.Range("K21").Value = 1500
$ws.Range("M21").Value = -1327

$ws.Range("H116").Value = 728.5
$ws.Range("I116").Value = 728.5
$ws.Range("K116").Value = 2185.5
$ws.Range("M116").Value = 1256.5

$ws.Range("H134").Value = 2343.5715
$ws.Range("I134").Value = 2343.5715
$ws.Range("K134").Value = 7030.7145
$ws.Range("M134").Value = -1960.7145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 100000
$ws.Range("J95").Value = 100000
$ws.Range("L95").Value = 100000
$ws.Range("N95").Value = -105492

$ws.Range("H131").Value = 102325.336
$ws.Range("J131").Value = 102325.336
$ws.Range("L131").Value = 102325.336
$ws.Range("N131").Value = -112405.336

$ws.Range("H132").Value = 3001.8823
$ws.Range("I132").Value = 2627
$ws.Range("J132").Value = 3901.6
$ws.Range("K132").Value = 7881
$ws.Range("L132").Value = 11704.8
$ws.Range("M132").Value = -5351
$ws.Range("N132").Value = -16764.8

$ws.Range("H136").Value = 63473.668
$ws.Range("J136").Value = 63473.668
$ws.Range("L136").Value = 190421.004
$ws.Range("N136").Value = -195521.004

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1144
$ws.Range("I55").Value = 1568.4286
$ws.Range("J55").Value = 772.625
$ws.Range("K55").Value = 1568.4286
$ws.Range("L55").Value = 772.625
$ws.Range("M55").Value = -1395.4286
$ws.Range("N55").Value = -1118.625

$ws.Range("H61").Value = 3846.1875
$ws.Range("I61").Value = 1349.3077
$ws.Range("K61").Value = 1349.3077
$ws.Range("M61").Value = -1147.3077

$ws.Range("H100").Value = 2999.3333
$ws.Range("J100").Value = 2999.5
$ws.Range("L100").Value = 2999.5
$ws.Range("N100").Value = -4081.5

$ws.Range("H113").Value = 3846.1875
$ws.Range("I113").Value = 1349.3077
$ws.Range("K113").Value = 1349.3077
$ws.Range("M113").Value = 820.6922999999999

$ws.Range("H122").Value = 8399.6
$ws.Range("I122").Value = 5334.3335
$ws.Range("K122").Value = 16003.0005
$ws.Range("M122").Value = -13553.0005

$ws.Range("H132").Value = 10579.333
$ws.Range("I132").Value = 12056.75
$ws.Range("K132").Value = 36170.25
$ws.Range("M132").Value = -33640.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 15007.223
$ws.Range("I51").Value = 19213
$ws.Range("K51").Value = 19213
$ws.Range("M51").Value = -18703

$ws.Range("H81").Value = 10066.477
$ws.Range("I81").Value = 10199.625
$ws.Range("K81").Value = 20399.25
$ws.Range("M81").Value = -19338.25

$ws.Range("H84").Value = 10066.477
$ws.Range("I84").Value = 10199.625
$ws.Range("K84").Value = 101996.25
$ws.Range("M84").Value = -96692.25

$ws.Range("H100").Value = 598.5454999999999
$ws.Range("I100").Value = 598.5454999999999
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1197.091
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -656.0909999999999
$ws.Range("N100").Value = ""

$ws.Range("H107").Value = 1836.909
$ws.Range("J107").Value = 3229.8
$ws.Range("L107").Value = 9689.400000000001
$ws.Range("N107").Value = -13529.4

$ws.Range("H113").Value = 544.2222
$ws.Range("I113").Value = 542.8570999999999
$ws.Range("J113").Value = 549
$ws.Range("K113").Value = 1628.5713
$ws.Range("L113").Value = 1647
$ws.Range("M113").Value = 541.4287000000002
$ws.Range("N113").Value = -5987

$ws.Range("H136").Value = 35367.87
$ws.Range("I136").Value = 2766.4285
$ws.Range("K136").Value = 8299.2855
$ws.Range("M136").Value = -5749.2855

$ws.Range("H137").Value = 149990
$ws.Range("J137").Value = 149990
$ws.Range("L137").Value = 149990
$ws.Range("N137").Value = -160190
